$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new column before column D ("Supplier"), shifting the rest right ---
$ws.Columns.Item(4).EntireColumn.Insert()

# The newly inserted column D should keep the same width as column C (both 13 chars)
$ws.Columns.Item(4).ColumnWidth = $ws.Columns.Item(3).ColumnWidth

# --- Give the new column a header: "Buyer" ---
$ws.Cells.Item(7, 4).Value = "Buyer"

# --- Fix up the workbook-level defined names so they again point at the right columns ---
# (Date / Inspection / Warehouse stay in place because they were left of the inserted column)
$wb.Names.Item("Avg_bag_weight").RefersTo = "=Sheet1!`$L`$8:`$L`$1048576"
$wb.Names.Item("Bags").RefersTo          = "=Sheet1!`$I`$8:`$I`$1048576"
$wb.Names.Item("Count").RefersTo         = "=Sheet1!`$N`$8:`$N`$1048576"
$wb.Names.Item("Fiche").RefersTo         = "=Sheet1!`$F`$8:`$F`$1048576"
$wb.Names.Item("Kg").RefersTo            = "=Sheet1!`$J`$8:`$J`$1048576"
$wb.Names.Item("KOR").RefersTo           = "=Sheet1!`$M`$8:`$M`$1048576"
$wb.Names.Item("Moisture").RefersTo      = "=Sheet1!`$O`$8:`$O`$1048576"
$wb.Names.Item("Price").RefersTo         = "=Sheet1!`$H`$8:`$H`$1048576"
$wb.Names.Item("Rejects").RefersTo       = "=Sheet1!`$P`$8:`$P`$1048576"
$wb.Names.Item("Supplier").RefersTo      = "=Sheet1!`$E`$8:`$E`$1048576"
$wb.Names.Item("Truck").RefersTo         = "=Sheet1!`$G`$8:`$G`$1048576"
$wb.Names.Item("Value_CFA").RefersTo     = "=Sheet1!`$K`$8:`$K`$1048576"

# --- Re-apply the AutoFilter / hidden _FilterDatabase name over the widened header row ---
$ws.AutoFilterMode = $False
[void]$ws.Range("A7:P7").AutoFilter()
$ws.Names.Item("Sheet1!_FilterDatabase").RefersTo = "=Sheet1!`$A`$7:`$P`$7"

# --- Restore the selected cell (points at the new Supplier column position) ---
[void]$ws.Range("E8").Select()
